# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets,
# mirroring the source data refresh captured in the commit.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 2227
    $ws.Range("F3").Value = 1695
    $ws.Range("F5").Value = 1080
    $ws.Range("F6").Value = 763
    $ws.Range("F8").Value = 5803
}
